# Update example import files with new data
# The "Strasse" (C), "PLZ" (D) and "Ort" (E) columns are re-shuffled among the
# existing patient rows (2-26) on the "patients" sheet; all other columns are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  Strasse = "Ringstraße 10";          Plz = 51645; Ort = "Gummersbach" }
    @{ Row = 3;  Strasse = "Wiesenweg 4";             Plz = 51643; Ort = "Gummersbach" }
    @{ Row = 4;  Strasse = "Ringstraße 8";            Plz = 51645; Ort = "Gummersbach" }
    @{ Row = 5;  Strasse = "Birkenweg 10";            Plz = 51674; Ort = "Wiehl" }
    @{ Row = 6;  Strasse = "Ulmenstraße 19";          Plz = 51702; Ort = "Bergneustadt" }
    @{ Row = 7;  Strasse = "Eschenallee 4";           Plz = 51688; Ort = "Wipperfürth" }
    @{ Row = 8;  Strasse = "Buchenweg 3";             Plz = 51674; Ort = "Wiehl" }
    @{ Row = 9;  Strasse = "Kastanienweg 5";          Plz = 51643; Ort = "Gummersbach" }
    @{ Row = 10; Strasse = "Kiefernweg 5";            Plz = 51688; Ort = "Wipperfürth" }
    @{ Row = 11; Strasse = "Lärchenweg 1";            Plz = 51674; Ort = "Wiehl" }
    @{ Row = 12; Strasse = "Hindenburgstraße 6";      Plz = 51643; Ort = "Gummersbach" }
    @{ Row = 13; Strasse = "Moltkestraße 10";         Plz = 51643; Ort = "Gummersbach" }
    @{ Row = 14; Strasse = "Breiter Weg 19";          Plz = 51647; Ort = "Bergneustadt" }
    @{ Row = 15; Strasse = "Lüdenscheider Straße 8";  Plz = 51688; Ort = "Wipperfürth" }
    @{ Row = 16; Strasse = "Weiherplatz 7";           Plz = 51674; Ort = "Wiehl" }
    @{ Row = 17; Strasse = "Kölner Straße 33";        Plz = 51647; Ort = "Bergneustadt" }
    @{ Row = 18; Strasse = "Marktstraße 5";           Plz = 51688; Ort = "Wipperfürth" }
    @{ Row = 19; Strasse = "Homburger Straße 2";      Plz = 51674; Ort = "Wiehl" }
    @{ Row = 20; Strasse = "Eichenweg 15";            Plz = 51647; Ort = "Bergneustadt" }
    @{ Row = 21; Strasse = "Tannenstraße 28";         Plz = 51688; Ort = "Wipperfürth" }
    @{ Row = 22; Strasse = "Ahornweg 16";             Plz = 51643; Ort = "Gummersbach" }
    @{ Row = 23; Strasse = "Erlenweg 36";             Plz = 51645; Ort = "Gummersbach" }
    @{ Row = 24; Strasse = "Pappelstraße 25";         Plz = 51647; Ort = "Bergneustadt" }
    @{ Row = 25; Strasse = "Zedernweg 38";            Plz = 51643; Ort = "Gummersbach" }
    @{ Row = 26; Strasse = "Fichtenstraße 21";        Plz = 51702; Ort = "Bergneustadt" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = $entry.Strasse
    $ws.Cells.Item($r, 4).Value = $entry.Plz
    $ws.Cells.Item($r, 5).Value = $entry.Ort
}

# Restore the active cell/selection recorded in the saved view state.
$ws.Range("J35").Select()
